# Update Fbln1-Itgb1 LR-pair sheet with recomputed TPM-based NATMI metrics.
# Adds a new "MuSCs" sending-cluster block (rows 17-21) and refreshes every
# numeric column for the existing "ECs"/"FAPs"/"Inflammatory-Mac" blocks so the
# sheet now spans A1:T21 (previously A1:T16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs (sending) -> ECs (target)
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fbln1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.126711
$ws.Range("H2").Value = 3.380133
$ws.Range("I2").Value = 0.02794157075411019
$ws.Range("J2").Value = 0.02794157075411019
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 159.4836373333333
$ws.Range("N2").Value = 478.450912
$ws.Range("O2").Value = 0.2983285084902258
$ws.Range("P2").Value = 0.2983285084902258
$ws.Range("Q2").Value = 179.6919685034773
$ws.Range("R2").Value = 1617.227716531296
$ws.Range("S2").Value = 0.008335767127947806
$ws.Range("T2").Value = 0.008335767127947806

# Row 3: ECs (sending) -> FAPs (target)
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fbln1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.126711
$ws.Range("H3").Value = 3.380133
$ws.Range("I3").Value = 0.02794157075411019
$ws.Range("J3").Value = 0.02794157075411019
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 172.558497
$ws.Range("N3").Value = 517.675491
$ws.Range("O3").Value = 0.3227862111630279
$ws.Range("P3").Value = 0.3227862111630279
$ws.Range("Q3").Value = 194.423556713367
$ws.Range("R3").Value = 1749.812010420303
$ws.Range("S3").Value = 0.009019153757662896
$ws.Range("T3").Value = 0.009019153757662896

# Row 4: ECs (sending) -> Inflammatory-Mac (target)
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fbln1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.126711
$ws.Range("H4").Value = 3.380133
$ws.Range("I4").Value = 0.02794157075411019
$ws.Range("J4").Value = 0.02794157075411019
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 74.38770566666666
$ws.Range("N4").Value = 223.163117
$ws.Range("O4").Value = 0.1391489036280481
$ws.Range("P4").Value = 0.1391489036280482
$ws.Range("Q4").Value = 83.81344623939566
$ws.Range("R4").Value = 754.3210161545611
$ws.Range("S4").Value = 0.003888038936079966
$ws.Range("T4").Value = 0.003888038936079967

# Row 5: ECs (sending) -> MuSCs (target)
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Fbln1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.126711
$ws.Range("H5").Value = 3.380133
$ws.Range("I5").Value = 0.02794157075411019
$ws.Range("J5").Value = 0.02794157075411019
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 58.41461433333333
$ws.Range("N5").Value = 175.243843
$ws.Range("O5").Value = 0.1092697975759847
$ws.Range("P5").Value = 0.1092697975759848
$ws.Range("Q5").Value = 65.81638853012433
$ws.Range("R5").Value = 592.347496771119
$ws.Range("S5").Value = 0.003053169780256676
$ws.Range("T5").Value = 0.003053169780256676

# Row 6: ECs (sending) -> Resolving-Mac (target)
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Fbln1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.126711
$ws.Range("H6").Value = 3.380133
$ws.Range("I6").Value = 0.02794157075411019
$ws.Range("J6").Value = 0.02794157075411019
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 69.746216
$ws.Range("N6").Value = 209.238648
$ws.Range("O6").Value = 0.1304665791427133
$ws.Range("P6").Value = 0.1304665791427133
$ws.Range("Q6").Value = 78.58382877557601
$ws.Range("R6").Value = 707.2544589801842
$ws.Range("S6").Value = 0.00364544115216284
$ws.Range("T6").Value = 0.003645441152162841

# Row 7: FAPs (sending) -> ECs (target)
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fbln1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 37.08140566666666
$ws.Range("H7").Value = 111.244217
$ws.Range("I7").Value = 0.9195904895727732
$ws.Range("J7").Value = 0.9195904895727732
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 159.4836373333333
$ws.Range("N7").Value = 478.450912
$ws.Range("O7").Value = 0.2983285084902258
$ws.Range("P7").Value = 0.2983285084902258
$ws.Range("Q7").Value = 5913.877453152878
$ws.Range("R7").Value = 53224.8970783759
$ws.Range("S7").Value = 0.274340059176042
$ws.Range("T7").Value = 0.274340059176042

# Row 8: FAPs (sending) -> FAPs (target)
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Fbln1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 37.08140566666666
$ws.Range("H8").Value = 111.244217
$ws.Range("I8").Value = 0.9195904895727732
$ws.Range("J8").Value = 0.9195904895727732
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 172.558497
$ws.Range("N8").Value = 517.675491
$ws.Range("O8").Value = 0.3227862111630279
$ws.Range("P8").Value = 0.3227862111630279
$ws.Range("Q8").Value = 6398.711628487282
$ws.Range("R8").Value = 57588.40465638554
$ws.Range("S8").Value = 0.2968311299507494
$ws.Range("T8").Value = 0.2968311299507494

# Row 9: FAPs (sending) -> Inflammatory-Mac (target)
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Fbln1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 37.08140566666666
$ws.Range("H9").Value = 111.244217
$ws.Range("I9").Value = 0.9195904895727732
$ws.Range("J9").Value = 0.9195904895727732
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 74.38770566666666
$ws.Range("N9").Value = 223.163117
$ws.Range("O9").Value = 0.1391489036280481
$ws.Range("P9").Value = 0.1391489036280482
$ws.Range("Q9").Value = 2758.400690438265
$ws.Range("R9").Value = 24825.60621394439
$ws.Range("S9").Value = 0.1279600084108314
$ws.Range("T9").Value = 0.1279600084108314

# Row 10: FAPs (sending) -> MuSCs (target)
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Fbln1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 37.08140566666666
$ws.Range("H10").Value = 111.244217
$ws.Range("I10").Value = 0.9195904895727732
$ws.Range("J10").Value = 0.9195904895727732
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 58.41461433333333
$ws.Range("N10").Value = 175.243843
$ws.Range("O10").Value = 0.1092697975759847
$ws.Range("P10").Value = 0.1092697975759848
$ws.Range("Q10").Value = 2166.096010956214
$ws.Range("R10").Value = 19494.86409860593
$ws.Range("S10").Value = 0.1004834666484176
$ws.Range("T10").Value = 0.1004834666484177

# Row 11: FAPs (sending) -> Resolving-Mac (target)
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Fbln1"
$ws.Range("C11").Value = "Itgb1"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 37.08140566666666
$ws.Range("H11").Value = 111.244217
$ws.Range("I11").Value = 0.9195904895727732
$ws.Range("J11").Value = 0.9195904895727732
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 69.746216
$ws.Range("N11").Value = 209.238648
$ws.Range("O11").Value = 0.1304665791427133
$ws.Range("P11").Value = 0.1304665791427133
$ws.Range("Q11").Value = 2586.287729210957
$ws.Range("R11").Value = 23276.58956289861
$ws.Range("S11").Value = 0.1199758253867327
$ws.Range("T11").Value = 0.1199758253867327

# Row 12: Inflammatory-Mac (sending) -> ECs (target)
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Fbln1"
$ws.Range("C12").Value = "Itgb1"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.02672533333333333
$ws.Range("H12").Value = 0.080176
$ws.Range("I12").Value = 0.000662767819130649
$ws.Range("J12").Value = 0.000662767819130649
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 159.4836373333333
$ws.Range("N12").Value = 478.450912
$ws.Range("O12").Value = 0.2983285084902258
$ws.Range("P12").Value = 0.2983285084902258
$ws.Range("Q12").Value = 4.262253368945778
$ws.Range("R12").Value = 38.360280320512
$ws.Range("S12").Value = 0.0001977225349565663
$ws.Range("T12").Value = 0.0001977225349565663

# Row 13: Inflammatory-Mac (sending) -> FAPs (target)
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Fbln1"
$ws.Range("C13").Value = "Itgb1"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.02672533333333333
$ws.Range("H13").Value = 0.080176
$ws.Range("I13").Value = 0.000662767819130649
$ws.Range("J13").Value = 0.000662767819130649
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 172.558497
$ws.Range("N13").Value = 517.675491
$ws.Range("O13").Value = 0.3227862111630279
$ws.Range("P13").Value = 0.3227862111630279
$ws.Range("Q13").Value = 4.611683351823999
$ws.Range("R13").Value = 41.50515016641599
$ws.Range("S13").Value = 0.0002139323132179652
$ws.Range("T13").Value = 0.0002139323132179652

# Row 14: Inflammatory-Mac (sending) -> Inflammatory-Mac (target)
$ws.Range("A14").Value = "Inflammatory-Mac"
$ws.Range("B14").Value = "Fbln1"
$ws.Range("C14").Value = "Itgb1"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.02672533333333333
$ws.Range("H14").Value = 0.080176
$ws.Range("I14").Value = 0.000662767819130649
$ws.Range("J14").Value = 0.000662767819130649
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 74.38770566666666
$ws.Range("N14").Value = 223.163117
$ws.Range("O14").Value = 0.1391489036280481
$ws.Range("P14").Value = 0.1391489036280482
$ws.Range("Q14").Value = 1.988036229843555
$ws.Range("R14").Value = 17.892326068592
$ws.Range("S14").Value = 0.00009222341539198231
$ws.Range("T14").Value = 0.00009222341539198234

# Row 15: Inflammatory-Mac (sending) -> MuSCs (target)
$ws.Range("A15").Value = "Inflammatory-Mac"
$ws.Range("B15").Value = "Fbln1"
$ws.Range("C15").Value = "Itgb1"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.02672533333333333
$ws.Range("H15").Value = 0.080176
$ws.Range("I15").Value = 0.000662767819130649
$ws.Range("J15").Value = 0.000662767819130649
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 58.41461433333333
$ws.Range("N15").Value = 175.243843
$ws.Range("O15").Value = 0.1092697975759847
$ws.Range("P15").Value = 0.1092697975759848
$ws.Range("Q15").Value = 1.561150039596444
$ws.Range("R15").Value = 14.050350356368
$ws.Range("S15").Value = 0.00007242050543628289
$ws.Range("T15").Value = 0.0000724205054362829

# Row 16: Inflammatory-Mac (sending) -> Resolving-Mac (target)
$ws.Range("A16").Value = "Inflammatory-Mac"
$ws.Range("B16").Value = "Fbln1"
$ws.Range("C16").Value = "Itgb1"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.02672533333333333
$ws.Range("H16").Value = 0.080176
$ws.Range("I16").Value = 0.000662767819130649
$ws.Range("J16").Value = 0.000662767819130649
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 69.746216
$ws.Range("N16").Value = 209.238648
$ws.Range("O16").Value = 0.1304665791427133
$ws.Range("P16").Value = 0.1304665791427133
$ws.Range("Q16").Value = 1.863990871338667
$ws.Range("R16").Value = 16.775917842048
$ws.Range("S16").Value = 0.00008646905012785232
$ws.Range("T16").Value = 0.00008646905012785234

# Row 17: MuSCs (sending) -> ECs (target)
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Fbln1"
$ws.Range("C17").Value = "Itgb1"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 2.088982666666667
$ws.Range("H17").Value = 6.266948
$ws.Range("I17").Value = 0.05180517185398602
$ws.Range("J17").Value = 0.05180517185398601
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 159.4836373333333
$ws.Range("N17").Value = 478.450912
$ws.Range("O17").Value = 0.2983285084902258
$ws.Range("P17").Value = 0.2983285084902258
$ws.Range("Q17").Value = 333.1585540062863
$ws.Range("R17").Value = 2998.426986056576
$ws.Range("S17").Value = 0.01545495965127947
$ws.Range("T17").Value = 0.01545495965127947

# Row 18: MuSCs (sending) -> FAPs (target)
$ws.Range("A18").Value = "MuSCs"
$ws.Range("B18").Value = "Fbln1"
$ws.Range("C18").Value = "Itgb1"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 2.088982666666667
$ws.Range("H18").Value = 6.266948
$ws.Range("I18").Value = 0.05180517185398602
$ws.Range("J18").Value = 0.05180517185398601
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 172.558497
$ws.Range("N18").Value = 517.675491
$ws.Range("O18").Value = 0.3227862111630279
$ws.Range("P18").Value = 0.3227862111630279
$ws.Range("Q18").Value = 360.471709219052
$ws.Range("R18").Value = 3244.245382971468
$ws.Range("S18").Value = 0.01672199514139768
$ws.Range("T18").Value = 0.01672199514139768

# Row 19: MuSCs (sending) -> Inflammatory-Mac (target)
$ws.Range("A19").Value = "MuSCs"
$ws.Range("B19").Value = "Fbln1"
$ws.Range("C19").Value = "Itgb1"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 2.088982666666667
$ws.Range("H19").Value = 6.266948
$ws.Range("I19").Value = 0.05180517185398602
$ws.Range("J19").Value = 0.05180517185398601
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 74.38770566666666
$ws.Range("N19").Value = 223.163117
$ws.Range("O19").Value = 0.1391489036280481
$ws.Range("P19").Value = 0.1391489036280482
$ws.Range("Q19").Value = 155.3946277507684
$ws.Range("R19").Value = 1398.551649756916
$ws.Range("S19").Value = 0.007208632865744772
$ws.Range("T19").Value = 0.007208632865744772

# Row 20: MuSCs (sending) -> MuSCs (target)
$ws.Range("A20").Value = "MuSCs"
$ws.Range("B20").Value = "Fbln1"
$ws.Range("C20").Value = "Itgb1"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 2.088982666666667
$ws.Range("H20").Value = 6.266948
$ws.Range("I20").Value = 0.05180517185398602
$ws.Range("J20").Value = 0.05180517185398601
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 58.41461433333333
$ws.Range("N20").Value = 175.243843
$ws.Range("O20").Value = 0.1092697975759847
$ws.Range("P20").Value = 0.1092697975759848
$ws.Range("Q20").Value = 122.0271168223516
$ws.Range("R20").Value = 1098.244051401164
$ws.Range("S20").Value = 0.005660740641874155
$ws.Range("T20").Value = 0.005660740641874155

# Row 21: MuSCs (sending) -> Resolving-Mac (target)
$ws.Range("A21").Value = "MuSCs"
$ws.Range("B21").Value = "Fbln1"
$ws.Range("C21").Value = "Itgb1"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 3
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 2.088982666666667
$ws.Range("H21").Value = 6.266948
$ws.Range("I21").Value = 0.05180517185398602
$ws.Range("J21").Value = 0.05180517185398601
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 69.746216
$ws.Range("N21").Value = 209.238648
$ws.Range("O21").Value = 0.1304665791427133
$ws.Range("P21").Value = 0.1304665791427133
$ws.Range("Q21").Value = 145.6986362895894
$ws.Range("R21").Value = 1311.287726606304
$ws.Range("S21").Value = 0.006758843553689931
$ws.Range("T21").Value = 0.006758843553689931
